$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated data values (2019 / 2020 rows revised) ---
$ws.Range("B9").Value = 3.3597913125480217
$ws.Range("C9").Value = 0.92105447730674861
$ws.Range("D9").Value = 2.1786611303191705

$ws.Range("B10").Value = -2.9778862191424338
$ws.Range("C10").Value = -4.178975398285079
$ws.Range("D10").Value = -9.3626431058061232

# --- Column widths (best-fit to header/content, matching Excel's AutoFit result) ---
$ws.Columns.Item(1).ColumnWidth = 4.166666666666667
$ws.Columns.Item(2).ColumnWidth = 5
$ws.Columns.Item(3).ColumnWidth = 8.5
$ws.Columns.Item(4).ColumnWidth = 6.833333333333333

# --- Selection now spans the full data block below the header ---
$ws.Range("A2:D12").Select() | Out-Null
